# Insert two new data rows at 319-320 (rest of the table, rows 319..419,
# shifts down to 321..421). This matches the diff: every existing row from
# 319 downward is now the same as the row that used to be two positions
# above it, and two brand-new rows of data appear at the former row 319.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("319:320").Insert()

# New row 319
$ws.Range("A319").Value = 7
$ws.Range("B319").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C319").Value = "Ñuble"
$ws.Range("D319").Value = 44627
$ws.Range("E319").Value = 16
$ws.Range("F319").Value = 100112020
$ws.Range("G319").Value = "Tomate"
$ws.Range("H319").Value = "Larga vida"
$ws.Range("I319").Value = "Primera"
$ws.Range("J319").Value = 600
$ws.Range("K319").Value = 9000
$ws.Range("L319").Value = 10000
$ws.Range("M319").Value = 9500
$ws.Range("N319").Value = "$/bandeja 18 kilos"
$ws.Range("O319").Value = "Región del Maule"
$ws.Range("P319").Value = 528
$ws.Range("Q319").Value = 18
$ws.Range("R319").Value = "Hortaliza"

# New row 320
$ws.Range("A320").Value = 7
$ws.Range("B320").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C320").Value = "Ñuble"
$ws.Range("D320").Value = 44627
$ws.Range("E320").Value = 16
$ws.Range("F320").Value = 100112020
$ws.Range("G320").Value = "Tomate"
$ws.Range("H320").Value = "Larga vida"
$ws.Range("I320").Value = "Segunda"
$ws.Range("J320").Value = 200
$ws.Range("K320").Value = 8000
$ws.Range("L320").Value = 8000
$ws.Range("M320").Value = 8000
$ws.Range("N320").Value = "$/bandeja 18 kilos"
$ws.Range("O320").Value = "Región del Maule"
$ws.Range("P320").Value = 444
$ws.Range("Q320").Value = 18
$ws.Range("R320").Value = "Hortaliza"
